$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (Bank Account), shifting existing
# columns F:L to G:M. This makes room for the new "Residency" column.
$ws.Columns.Item(6).Insert()

# Match the width of the neighboring "Address"/"KYC Type" columns (11)
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Rename the existing "Type" header (column E) to "KYC Type"
$ws.Range("E1").Value = "KYC Type"

# Populate the newly inserted "Residency" column (F)
$ws.Range("F1").Value = "Residency"
$ws.Range("F2").Value = "Domestic"
$ws.Range("F3").Value = "Foreign"

# Match the new active cell selection from the diff
$ws.Range("F4").Select()
